$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name order in column A (shared-string reorder cascade) ---
$ws.Range("A63").Value = "Hungria"
$ws.Range("A64").Value = "Nueva Zelanda"
$ws.Range("A65").Value = "Irak"
$ws.Range("A148").Value = "Gabon"
$ws.Range("A149").Value = "Liberia"
$ws.Range("A168").Value = "Maldivas"
$ws.Range("A169").Value = "Sudan"
$ws.Range("A170").Value = "Angola"
$ws.Range("A210").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A211").Value = "Anguila"
$ws.Range("A212").Value = "Islas Virgenes Britanicas"

# --- Update last-updated timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 07:52"

# --- Update numeric data cells ---
$ws.Range("B4").Value = 533115
$ws.Range("C4").Value = 236
$ws.Range("E4").Value = 482033
$ws.Range("B25").Value = 8504
$ws.Range("C25").Value = 58
$ws.Range("D25").Value = 972
$ws.Range("E25").Value = 7243
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 289
$ws.Range("B63").Value = 1410
$ws.Range("C63").Value = 100
$ws.Range("D63").Value = 118
$ws.Range("E63").Value = 1207
$ws.Range("F63").Value = 17
$ws.Range("H63").Value = 85
$ws.Range("B64").Value = 1330
$ws.Range("C64").Value = 18
$ws.Range("D64").Value = 471
$ws.Range("E64").Value = 855
$ws.Range("F64").Value = 5
$ws.Range("H64").Value = 4
$ws.Range("B65").Value = 1318
$ws.Range("D65").Value = 601
$ws.Range("E65").Value = 645
$ws.Range("F65").Value = 0
$ws.Range("H65").Value = 72
$ws.Range("B148").Value = 49
$ws.Range("C148").Value = 3
$ws.Range("D148").Value = 1
$ws.Range("E148").Value = 47
$ws.Range("H148").Value = 1
$ws.Range("B149").Value = 48
$ws.Range("D149").Value = 3
$ws.Range("E149").Value = 40
$ws.Range("H149").Value = 5
$ws.Range("B168").Value = 20
$ws.Range("C168").Value = 1
$ws.Range("D168").Value = 13
$ws.Range("E168").Value = 7
$ws.Range("H168").Value = 0
$ws.Range("D169").Value = 2
$ws.Range("E169").Value = 15
$ws.Range("D170").Value = 4
$ws.Range("E170").Value = 13
$ws.Range("H170").Value = 2
$ws.Range("C210").Value = 1
$ws.Range("D211").Value = 0
$ws.Range("E211").Value = 3
$ws.Range("B212").Value = 3
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 2
$ws.Range("E212").Value = 1
